# contraordenarcheque.xlsx - "cambio mapeo para ios actualizar datos de seguridad"
# Update the security/cheque-range data used by the datadriven test so the
# reserved cheque numbers no longer collide with ones already consumed by
# another (iOS) test run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Row 2 -> numeroCheque (column P)
$ws.Range("P2").Value = 65728

# Row 3 -> rangoDesde / rangoHasta (columns Q / R)
$ws.Range("Q3").Value = 65729
$ws.Range("R3").Value = 65730

# Reflect the editor's on-screen state at the moment the file was saved:
# scrolled one column further left and the cursor left on the last data cell.
[void]$ws.Activate()
[void]$ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("Q4").Select()
